# Update the cryptocurrency price/volume table to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.196.44"
$ws.Range("E2").Value = "  +10.67%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.479.62"
$ws.Range("E3").Value = "  +6.87%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "416.96"
$ws.Range("E5").Value = "  +4.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "124.77"
$ws.Range("E6").Value = "  +15.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.470.07"
$ws.Range("E7").Value = "  +6.80%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.598"
$ws.Range("E8").Value = "  +3.76%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.679"
$ws.Range("E10").Value = "  +9.88%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.130"
$ws.Range("E11").Value = "  +36.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "41.78"
$ws.Range("E12").Value = "  +6.62%  "

$ws.Range("E13").Value = "  -0.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.018.69"
$ws.Range("E14").Value = "  +6.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.60"
$ws.Range("E15").Value = "  +4.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.91"
$ws.Range("E16").Value = "  +5.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.468.51"
$ws.Range("E17").Value = "  +6.71%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.985.17"
$ws.Range("E18").Value = "  +10.66%  "

$ws.Range("E19").Value = "  +0.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.20"
$ws.Range("E20").Value = "  +2.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000135"
$ws.Range("E21").Value = "  +26.06%  "

$ws.Range("E22").Value = "  +1.72%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "319.71"
$ws.Range("E23").Value = "  +8.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.01"
$ws.Range("E24").Value = "  +12.44%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.02"
$ws.Range("E25").Value = "  +0.87%  "

$ws.Range("E26").Value = "  +1.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "31.23"
$ws.Range("E27").Value = "  +11.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.83"
$ws.Range("E28").Value = "  +5.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.93"
$ws.Range("E29").Value = "  +0.92%  "

$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.175"
$ws.Range("E30").Value = "  +3.78%  "

$ws.Range("B31").Value = "LEO"
$ws.Range("C31").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.31"
$ws.Range("E31").Value = "  -1.64%  "

$ws.Range("E32").Value = "  +3.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.64"
$ws.Range("E33").Value = "  +4.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.47"
$ws.Range("E34").Value = "  +6.19%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.57"
$ws.Range("E35").Value = "  +19.89%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.01"
$ws.Range("E36").Value = "  +0.73%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0491"
$ws.Range("E37").Value = "  -0.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.13"
$ws.Range("E38").Value = "  +1.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.50"
$ws.Range("E39").Value = "  +1.27%  "

$ws.Range("E40").Value = "  -0.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.07"
$ws.Range("E41").Value = "  +2.77%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.02"
$ws.Range("E42").Value = "  +8.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.127"
$ws.Range("E43").Value = "  +4.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "135.05"
$ws.Range("E44").Value = "  -1.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.40"
$ws.Range("E45").Value = "  +4.39%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.286"
$ws.Range("E46").Value = "  +1.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.94"
$ws.Range("E47").Value = "  +1.18%  "

$ws.Range("E48").Value = "  +1.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.23"
$ws.Range("E49").Value = "  -0.29%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.816.23"
$ws.Range("E50").Value = "  +6.60%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.183.31"
$ws.Range("E51").Value = "  +2.15%  "
